# Applies the commit: "add column to game_nfl, add one week of test schedule"
#  1. Adds a new worksheet "test_schedule" after "team_nfl"
#  2. Populates it with headers + one week (16 games) of 2016 NFL schedule/test data
#  3. Turns the range into a table (Table3) styled like the existing team table
#  4. Formats the start_time column with a custom date/time number format

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. add the new sheet right after team_nfl -----------------------------
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "test_schedule"

# --- 2. headers --------------------------------------------------------------
$headers = @("Season", "week", "home_id", "away_id", "home_score", "away_score", "completed", "start_time")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 3. one week (week 1, 2016) of game rows ---------------------------------
# columns: Season, week, home_id, away_id, home_score(formula), away_score(formula), start_time(serial)
$games = @(
    @(2016, 1, 1,  2,  0,  38, 42621.833333333336),
    @(2016, 1, 3,  4,  49, 32, 42624.541666666664),
    @(2016, 1, 5,  6,  27, 18, 42624.541666666664),
    @(2016, 1, 7,  8,  42, 30, 42624.541666608799),
    @(2016, 1, 9,  10, 17, 29, 42624.541666608799),
    @(2016, 1, 11, 12, 15, 38, 42624.541666608799),
    @(2016, 1, 13, 14, 6,  5,  42624.541666608799),
    @(2016, 1, 15, 16, 4,  39, 42624.541666608799),
    @(2016, 1, 17, 18, 17, 6,  42624.541666608799),
    @(2016, 1, 19, 20, 8,  36, 42624.541666608799),
    @(2016, 1, 21, 22, 25, 11, 42624.541666608799),
    @(2016, 1, 23, 24, 26, 43, 42624.541666608799),
    @(2016, 1, 25, 26, 34, 40, 42624.666666666664),
    @(2016, 1, 27, 28, 41, 19, 42624.854166666664),
    @(2016, 1, 29, 30, 2,  9,  42625.833333333336),
    @(2016, 1, 31, 0,  42, 33, 42625.833333333336)
)

# number format for start_time, applied up front so the style created for
# the first date cell is reused (same xf) for every subsequent one
$dateFormat = "[`$-409]m/d/yy\ h:mm\ AM/PM;@"

$row = 2
foreach ($g in $games) {
    $ws.Cells.Item($row, 1).Value = $g[0]
    $ws.Cells.Item($row, 2).Value = $g[1]
    $ws.Cells.Item($row, 3).Value = $g[2]
    $ws.Cells.Item($row, 4).Value = $g[3]
    $ws.Cells.Item($row, 5).Formula = "=RANDBETWEEN(0, 50)"
    $ws.Cells.Item($row, 6).Formula = "=RANDBETWEEN(0, 50)"
    $ws.Cells.Item($row, 7).Value = $true
    $ws.Cells.Item($row, 8).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 8).Value = $g[6]
    $row++
}

$lastRow = $row - 1

# --- 4. column widths ---------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.83
$ws.Columns.Item(6).ColumnWidth = 12.33
$ws.Columns.Item(7).ColumnWidth = 11.5
$ws.Columns.Item(8).ColumnWidth = 19.67

# --- 5. turn the range into a table, matching the team_nfl table's style ----
$tableRange = $ws.Range("A1:H$lastRow")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table3"
$lo.TableStyle = "TableStyleLight9"

# --- 6. selection matching the final saved view of the author --------------
$ws.Range("N23").Select()
